$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows with new player/position/team data per the diff
$ws.Range("A2").Value = "Russell Westbrook"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Denver Nuggets"

$ws.Range("A4").Value = "Ben Simmons"
$ws.Range("B4").Value = "PG,C"
$ws.Range("C4").Value = "Brooklyn Nets"

$ws.Range("A5").Value = "Jalen Green"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Houston Rockets"

$ws.Range("A7").Value = "Pascal Siakam"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Indiana Pacers"

$ws.Range("A14").Value = "Jalen Suggs"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Orlando Magic"

$ws.Range("A15").Value = "Chris Paul"
$ws.Range("B15").Value = "PG"
$ws.Range("C15").Value = "San Antonio Spurs"

$ws.Range("A16").Value = "Jonathan Kuminga"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Golden State Warriors"
